# Hortaliza, Comercializadora del Agro de Limarí - Choclo.xlsx
# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record needs to be inserted into the data table,
# right before the existing row 97 (which shifts that row, and every row
# after it, down by one). The sheet's used range grows from A1:R131 to
# A1:R132 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 97; this pushes the former rows
# 97..131 down to 98..132 and preserves their formatting/styles.
$ws.Rows.Item(97).Insert()

# Fill the newly inserted row 97 with the new record's data.
$ws.Range("A97").Value2 = 2
$ws.Range("B97").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C97").Value2 = "Coquimbo"
$ws.Range("D97").Value2 = 44755
$ws.Range("E97").Value2 = 4
$ws.Range("F97").Value2 = 100112024
$ws.Range("G97").Value2 = "Choclo"
$ws.Range("H97").Value2 = "Dulce o Americano"
$ws.Range("I97").Value2 = "Primera"
$ws.Range("J97").Value2 = 400
$ws.Range("K97").Value2 = 29000
$ws.Range("L97").Value2 = 30000
$ws.Range("M97").Value2 = 29500
$ws.Range("N97").Value2 = "$/malla 60 unidades"
$ws.Range("O97").Value2 = "Provincia de Limarí"
$ws.Range("P97").Value2 = 492
$ws.Range("Q97").Value2 = 60
$ws.Range("R97").Value2 = "Hortaliza"
